# Generate Report for Handoff
# - Overview sheet: e2fbda3d-...md row status changes from
#   "Handed back: in sync with en-US" to "Ready for handoff"
# - zh-cn sheet: Latest Handoff Datetime for e2fbda3d row refreshed,
#   and (since both rows shared the same handoff time) the
#   2a54065f row's Latest Handoff Datetime also refreshes to the
#   same new timestamp; Status for e2fbda3d row becomes "Ready for handoff"
# - de-de sheet: same pattern with its own new timestamp

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the e2fbda3d-...md file ---
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row 2 (2a54065f) and row 3 (e2fbda3d) ---
$wsZhCn.Range("D2").Value = "2016-03-04 11:25:20"
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-03-04 11:25:20"

# --- de-de sheet: row 2 (2a54065f) and row 3 (e2fbda3d) ---
$wsDeDe.Range("D2").Value = "2016-03-04 11:25:33"
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-03-04 11:25:33"
